$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New shared strings / cell values
# ---------------------------------------------------------------------------
$ws.Range("B14").Value2 = 172

$ws.Range("A17").Value2 = "Total time"
$ws.Range("B17").Formula = "=SUM(B2:B16)"

$ws.Range("E18").Value2 = "2-Sigma Range"
$ws.Range("F18").Formula = "=F17-2*(G17^0.5)"
$ws.Range("G18").Formula = "=F17+2*(G17^0.5)"

# ---------------------------------------------------------------------------
# 2. Apply the "Output" cell style (builtin 21 / German "Ausgabe") to the
#    whole populated table, then overlay the 0.00 number format on the two
#    numeric (computed) columns.
# ---------------------------------------------------------------------------
$ws.Range("A1:G1").Style = "Output"
$ws.Range("A2:E17").Style = "Output"
$ws.Range("F2:G18").Style = "Output"
$ws.Range("F2:G18").NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 3. "Total time" / "2-Sigma Range" label cells get the bold-italic variant
#    of the Output style, right aligned / vertically centered.
#    (applied per contiguous area - union ranges don't reliably propagate
#    .Style across all their areas)
# ---------------------------------------------------------------------------
$labelAreas = @($ws.Range("A17"), $ws.Range("D18:E18"))
foreach ($area in $labelAreas) {
    $area.Style = "Output"
    $area.Font.Italic = $true
    $area.HorizontalAlignment = -4152
    $area.VerticalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.8333333333333
$ws.Columns.Item(3).ColumnWidth = 9.66666666666667
$ws.Columns.Item(4).ColumnWidth = 11.1666666666667
$ws.Columns.Item(5).ColumnWidth = 10.6666666666667
$ws.Columns.Item(6).ColumnWidth = 7.33333333333333
$ws.Columns.Item(7).ColumnWidth = 7.16666666666667

# ---------------------------------------------------------------------------
# 5. Row heights (all data rows grow from 12.75 to 15 once the bigger
#    11pt Calibri "Output" font is applied)
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 18; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

# ---------------------------------------------------------------------------
# 6. Selection / zoom
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("H6").Select()
$excel.ActiveWindow.Zoom = 175
